# sistemate personae e presentazione
# Slide 4 ("Risultati posizionamento"): the content placeholder had a
# stray "Immagini grafici" line followed by an empty paragraph. Remove
# the text paragraph, leaving a single empty paragraph with bullets
# turned off (matches PowerPoint's own behaviour when the last run of
# text is cleared from a placeholder).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Remove the "Immagini grafici" paragraph entirely (also removes the
# paragraph break, merging what remains into a single paragraph).
$firstPara = $tr.Paragraphs(1, 1)
$firstPara.Delete()

# The remaining (now only) paragraph keeps its bullet off, as
# PowerPoint does for an emptied placeholder paragraph.
$shape.TextFrame.TextRange.ParagraphFormat.Bullet.Visible = $false
